$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column B
$ws.Range("B1").Value = "Active Cases"

# Full replacement data set (Cluster Name, Active Cases) for rows 2..38
$data = @(
    @('3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North', 15),
    @('3376 Royal Freemasons Coppin Centre Melbourne', 10),
    @('3601 Baptcare Westhaven Community', 29),
    @('3653 Fronditha Thalpori St Albans Aged Care', 39),
    @('44121 Wallan Primary School Wallan', 18),
    @('44165 Greenvale Primary School', 38),
    @('44234 Lucknow Primary School Bairnsdale', 26),
    @('44395 Buln Buln Primary School', 11),
    @('44554 Korumburra Primary School', 10),
    @('44667 Beaumaris Primary School Beaumaris', 22),
    @('44701 Hampton Park Primary School Hampton Park', 10),
    @('44811 Dandenong North Primary School Dandenong', 27),
    @('44820 Bourchier Street Primary School Shepparton', 11),
    @('44853 St Albans North Primary School', 13),
    @('45158 Rowellyn Park Primary School Carrum Downs', 16),
    @('45249 Creekside K-9 College Caroline Springs', 16),
    @('45695 Sacred Heart Primary School Yarrawonga', 20),
    @('4574 Village Glen Aged Care Residences Mornington', 18),
    @('45812 St Mary''s Primary School Hampton', 11),
    @('45967 St Clement of Rome School Bulleen', 10),
    @('46037 Nazareth Catholic Primary School Grovedale', 21),
    @('46050 Our Lady''s Catholic Primary School Craigieburn', 42),
    @('46125 Our Lady of the Southern Cross Primary School Manor Lakes', 29),
    @('46140 Caulfield Grammar School Malvern Campus Glen Iris', 10),
    @('46276 Hillcrest Christian College Clyde North', 16),
    @('46322 Minaret College Officer Campus Officer', 12),
    @('46328 Ilim College Dallas Main Campus Dallas', 30),
    @('46328 Ilim College Dallas Primary Campus Inverloch', 13),
    @('46390 Al Siraat College Epping', 45),
    @('50681 Broadmeadows Special Developmental School Broadmeadows', 11),
    @('52380 Al Iman College Melton South', 18),
    @('52473 John Henry Primary School Pakenham', 15),
    @('52786 Hume Anglican Grammar Donnybrook Campus', 16),
    @('Adass Israel School Elsternwick', 10),
    @('Hamilton Country Music Festival Hamilton Golf Club Hamilton', 15),
    @('Melton Willows Melton', 11),
    @('Sirius College Shepparton Campus Shepparton', 11)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
